$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1Updates = @(
    @(39, 8, 1191.6666),
    @(39, 9, 100),
    @(39, 10, 1555.5555),
    @(39, 11, 300),
    @(39, 12, 4666.666499999999),
    @(39, 13, -4),
    @(39, 14, -5258.666499999999),
    @(40, 8, 2204.25),
    @(40, 10, 2183),
    @(40, 12, 2183),
    @(40, 14, -2533),
    @(64, 8, 3133.3333),
    @(64, 9, 3000),
    @(64, 11, 3000),
    @(64, 13, -2752),
    @(67, 8, 3133.3333),
    @(67, 9, 3000),
    @(67, 11, 3000),
    @(67, 13, -2142),
    @(98, 8, 2133.8147),
    @(98, 9, 1818.4348),
    @(98, 10, 3947.25),
    @(98, 11, 1818.4348),
    @(98, 12, 3947.25),
    @(98, 13, -320.4348),
    @(98, 14, -6943.25),
    @(112, 8, 1817.7759),
    @(112, 9, 1500),
    @(112, 10, 1829.125),
    @(112, 11, 4500),
    @(112, 12, 5487.375),
    @(112, 13, -3392),
    @(112, 14, -7703.375),
    @(122, 8, 2133.8147),
    @(122, 9, 1818.4348),
    @(122, 10, 3947.25),
    @(122, 11, 5455.3044),
    @(122, 12, 11841.75),
    @(122, 13, -3005.3044),
    @(122, 14, -16741.75),
    @(138, 8, 2674.9434),
    @(138, 9, 2832.4138),
    @(138, 10, 2484.6667),
    @(138, 11, 8497.241399999999),
    @(138, 12, 7454.000100000001),
    @(138, 13, -3357.241399999999),
    @(138, 14, -17734.0001)
)
foreach ($u in $ws1Updates) {
    $r = $u[0]; $c = $u[1]; $v = $u[2]
    if ($null -eq $v) {
        $ws1.Cells.Item($r, $c).ClearContents()
    } else {
        $ws1.Cells.Item($r, $c).Value = $v
    }
}

$ws2 = $wb.Worksheets.Item("ARM")
$ws2Updates = @(
    @(32, 8, 4363.6665),
    @(32, 9, 3342.9207),
    @(32, 10, 15081.5),
    @(32, 11, 3342.9207),
    @(32, 12, 15081.5),
    @(32, 13, -3055.9207),
    @(32, 14, -15655.5),
    @(44, 8, 30000),
    @(44, 10, 30000),
    @(44, 12, 30000),
    @(44, 14, -30976),
    @(55, 8, 20000),
    @(55, 10, 20000),
    @(55, 12, 20000),
    @(55, 14, -20630),
    @(61, 8, 5895.207),
    @(61, 9, 5733.25),
    @(61, 11, 5733.25),
    @(61, 13, -5521.25),
    @(74, 8, 841.7692),
    @(74, 9, 535.1429000000001),
    @(74, 10, 3524.75),
    @(74, 11, 535.1429000000001),
    @(74, 12, 3524.75),
    @(74, 13, 338.8570999999999),
    @(74, 14, -5272.75),
    @(77, 8, 841.7692),
    @(77, 9, 535.1429000000001),
    @(77, 10, 3524.75),
    @(77, 11, 2675.7145),
    @(77, 12, 17623.75),
    @(77, 13, 1692.2855),
    @(77, 14, -26359.75),
    @(80, 8, 48500),
    @(80, 9, 0),
    @(80, 10, 48500),
    @(80, 11, 0),
    @(80, 12, 48500),
    @(80, 13, $null),
    @(80, 14, -50496),
    @(83, 8, 48500),
    @(83, 9, 0),
    @(83, 10, 48500),
    @(83, 11, 0),
    @(83, 12, 145500),
    @(83, 13, $null),
    @(83, 14, -155484),
    @(97, 8, 1466.0714),
    @(97, 9, 1482.7),
    @(97, 10, 1424.5),
    @(97, 11, 1482.7),
    @(97, 12, 1424.5),
    @(97, 13, -986.7),
    @(97, 14, -2416.5),
    @(102, 8, 1149.9375),
    @(102, 9, 985),
    @(102, 11, 985),
    @(102, 13, 637),
    @(136, 8, 5895.207),
    @(136, 9, 5733.25),
    @(136, 11, 17199.75),
    @(136, 13, -14649.75),
    @(139, 8, 42738.2),
    @(139, 10, 42738.2),
    @(139, 12, 42738.2),
    @(139, 14, -53018.2)
)
foreach ($u in $ws2Updates) {
    $r = $u[0]; $c = $u[1]; $v = $u[2]
    if ($null -eq $v) {
        $ws2.Cells.Item($r, $c).ClearContents()
    } else {
        $ws2.Cells.Item($r, $c).Value = $v
    }
}

$ws3 = $wb.Worksheets.Item("BSM")
$ws3Updates = @(
    @(86, 8, 146502.58),
    @(86, 9, 3919.6667),
    @(86, 11, 3919.6667),
    @(86, 13, -2796.6667),
    @(89, 8, 146502.58),
    @(89, 9, 3919.6667),
    @(89, 11, 19598.3335),
    @(89, 13, -13982.3335),
    @(134, 8, 6550.0435),
    @(134, 9, 7047.095),
    @(134, 11, 21141.285),
    @(134, 13, -18606.285)
)
foreach ($u in $ws3Updates) {
    $r = $u[0]; $c = $u[1]; $v = $u[2]
    if ($null -eq $v) {
        $ws3.Cells.Item($r, $c).ClearContents()
    } else {
        $ws3.Cells.Item($r, $c).Value = $v
    }
}

$ws4 = $wb.Worksheets.Item("CRP")
$ws4Updates = @(
    @(31, 8, 2309.5293),
    @(31, 9, 2025.5834),
    @(31, 11, 2025.5834),
    @(31, 13, -1730.5834),
    @(34, 8, 2309.5293),
    @(34, 9, 2025.5834),
    @(34, 11, 2025.5834),
    @(34, 13, -1823.5834),
    @(105, 8, 846),
    @(105, 9, 826.0909),
    @(105, 11, 826.0909),
    @(105, 13, 920.9091),
    @(107, 8, 657.619),
    @(107, 9, 409.7647),
    @(107, 11, 409.7647),
    @(107, 13, 1510.2353)
)
foreach ($u in $ws4Updates) {
    $r = $u[0]; $c = $u[1]; $v = $u[2]
    if ($null -eq $v) {
        $ws4.Cells.Item($r, $c).ClearContents()
    } else {
        $ws4.Cells.Item($r, $c).Value = $v
    }
}

$ws5 = $wb.Worksheets.Item("CUL")
$ws5Updates = @(
    @(107, 8, 421.08334),
    @(107, 10, 385.3),
    @(107, 12, 1155.9),
    @(107, 14, -4995.9),
    @(131, 8, 21322.234),
    @(131, 10, 21937.152),
    @(131, 12, 65811.45599999999),
    @(131, 14, -75891.45599999999)
)
foreach ($u in $ws5Updates) {
    $r = $u[0]; $c = $u[1]; $v = $u[2]
    if ($null -eq $v) {
        $ws5.Cells.Item($r, $c).ClearContents()
    } else {
        $ws5.Cells.Item($r, $c).Value = $v
    }
}

$ws6 = $wb.Worksheets.Item("GSM")
$ws6Updates = @(
    @(80, 8, 2666.6667),
    @(80, 9, 2544.9092),
    @(80, 10, 4006),
    @(80, 11, 2544.9092),
    @(80, 12, 4006),
    @(80, 13, -1546.9092),
    @(80, 14, -6002),
    @(83, 8, 2666.6667),
    @(83, 9, 2544.9092),
    @(83, 10, 4006),
    @(83, 11, 12724.546),
    @(83, 12, 20030),
    @(83, 13, -7732.546),
    @(83, 14, -30014),
    @(113, 8, 773.0952),
    @(113, 9, 569.0769),
    @(113, 10, 1104.625),
    @(113, 11, 569.0769),
    @(113, 12, 1104.625),
    @(113, 13, 1600.9231),
    @(113, 14, -5444.625)
)
foreach ($u in $ws6Updates) {
    $r = $u[0]; $c = $u[1]; $v = $u[2]
    if ($null -eq $v) {
        $ws6.Cells.Item($r, $c).ClearContents()
    } else {
        $ws6.Cells.Item($r, $c).Value = $v
    }
}

$ws7 = $wb.Worksheets.Item("LTW")
$ws7Updates = @(
    @(46, 8, 1538.125),
    @(46, 9, 1095.0834),
    @(46, 11, 1095.0834),
    @(46, 13, -907.0834),
    @(61, 8, 2873.5386),
    @(61, 9, 2528.3333),
    @(61, 10, 3650.25),
    @(61, 11, 2528.3333),
    @(61, 12, 3650.25),
    @(61, 13, -2326.3333),
    @(61, 14, -4054.25),
    @(113, 8, 2873.5386),
    @(113, 9, 2528.3333),
    @(113, 10, 3650.25),
    @(113, 11, 2528.3333),
    @(113, 12, 3650.25),
    @(113, 13, -358.3332999999998),
    @(113, 14, -7990.25)
)
foreach ($u in $ws7Updates) {
    $r = $u[0]; $c = $u[1]; $v = $u[2]
    if ($null -eq $v) {
        $ws7.Cells.Item($r, $c).ClearContents()
    } else {
        $ws7.Cells.Item($r, $c).Value = $v
    }
}

$ws8 = $wb.Worksheets.Item("WVR")
$ws8Updates = @(
    @(126, 8, 6405.091),
    @(126, 9, 10169.083),
    @(126, 10, 1888.3),
    @(126, 11, 30507.249),
    @(126, 12, 5664.9),
    @(126, 13, -28037.249),
    @(126, 14, -10604.9)
)
foreach ($u in $ws8Updates) {
    $r = $u[0]; $c = $u[1]; $v = $u[2]
    if ($null -eq $v) {
        $ws8.Cells.Item($r, $c).ClearContents()
    } else {
        $ws8.Cells.Item($r, $c).Value = $v
    }
}
